# Weekly update: a new week of "Ajo" (garlic) price data was collected at
# Mercado Mayorista Lo Valledor de Santiago. The new record slots in right
# after the current newest "Primera"/Chino garlic entry (old row 671),
# pushing that row and everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 671; Excel shifts rows 671:711 down to 672:712.
$ws.Rows(671).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(671, 1).Value  = 6
$ws.Cells.Item(671, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(671, 3).Value  = "Metropolitana"
$ws.Cells.Item(671, 4).Value  = 44706
$ws.Cells.Item(671, 5).Value  = 13
$ws.Cells.Item(671, 6).Value  = 100112003
$ws.Cells.Item(671, 7).Value  = "Ajo"
$ws.Cells.Item(671, 8).Value  = "Chino"
$ws.Cells.Item(671, 9).Value  = "Primera"
$ws.Cells.Item(671, 10).Value = 1900
$ws.Cells.Item(671, 11).Value = 16000
$ws.Cells.Item(671, 12).Value = 17000
$ws.Cells.Item(671, 13).Value = 16579
$ws.Cells.Item(671, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(671, 15).Value = "China"
$ws.Cells.Item(671, 16).Value = 1658
$ws.Cells.Item(671, 17).Value = 10
$ws.Cells.Item(671, 18).Value = "Hortaliza"
